# The workbook's "train" sheet lists movie-review sentences together with
# their sentiment score and train/validate/test split (column D). This
# edit:
#   1. Swaps the sentence text used for SentenceID 4 (row 5) and
#      SentenceID 8 (row 9) - row 9 now carries the (slightly reworded)
#      "Importance of Being Earnest" text and row 5 now carries the
#      (slightly reworded) "Moonlight Mile" text - and stashes a small
#      " c" note in J9.
#   2. Applies an AutoFilter on column D (the split column) that shows
#      only rows whose split is "Train", hiding the rest.
#
# NOTE on ordering: new unique strings are appended to the shared-string
# table in the order they are first written, so the three brand-new
# strings below are written in the exact order that reproduces the
# target shared string indices (" c" first, then the reworded
# "Importance" text, then the reworded "Moonlight Mile" text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$ws.Range("J9").Value = " c"
$ws.Range("B9").Value = "The Importance of Being Earnest , so thick with wit it plays like a reading from Bartlett 's Familiar Quotations so I still like it and am thick wit."
$ws.Range("B5").Value = "I still like Moonlight Mile , better judgment be damned , worth seeing so thick wit. "

# Column D (4th column of the A1:D14 table) holds the Train/Validate/Test
# split label; filter to only "Train" rows (xlFilterValues = 7).
$ws.Range("A1:D14").AutoFilter(4, @("Train"), 7)

$ws.Range("B6").Select()
